# Add a `brand_name` column in front of the existing `search_keyword` column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing search_keyword values from column A over to column B.
$ws.Range("B1").Value = $ws.Range("A1").Text
$ws.Range("B2").Value = $ws.Range("A2").Text
$ws.Range("B3").Value = $ws.Range("A3").Text
$ws.Range("B4").Value = $ws.Range("A4").Text
$ws.Range("B5").Value = $ws.Range("A5").Text

# New column A: brand values for each search keyword row, then the header.
$ws.Range("A2").Value = "Optimum Nutrition"
$ws.Range("A3").Value = "Optimum Nutrition"
$ws.Range("A4").Value = "Nature made"
$ws.Range("A5").Value = "Optimum Nutrition"
$ws.Range("A1").Value = "brand_name"

# Keep the active selection where Excel left it after editing.
$ws.Range("A5").Select()
